# The presentation's theme (ppt/theme/theme1.xml, "Integral" / "Red Violet")
# is swapped with the secondary theme (ppt/theme/theme2.xml, "Office Theme" /
# "Office") that is already embedded in the package for the notes master.
# Both themes share an identical font scheme and format scheme - they only
# differ in their 12-slot color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) and display name. Re-point the active theme's color scheme to the
# "Office" palette so the deck visually matches the target "Office Theme".

$p = $ppt.ActivePresentation

# Target palette (the "Office" color scheme that currently lives in theme2.xml)
# in the same dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order used by
# ThemeColorScheme.Colors(1..12).
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgbVal = $r + ($g * 256) + ($b * 65536)

    $c = $tcs.Colors($i)
    $c.RGB = $rgbVal
}

Write-Output "Theme color scheme updated to Office palette."
